$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    4   = 16.452
    6   = 16.322
    7   = 16.48
    8   = 16.44
    16  = 16.618
    20  = 16.24
    21  = 16.482
    28  = 17.006
    29  = 16.941
    30  = 16.321
    32  = 17.4
    40  = 16.498
    46  = 16.858
    51  = 16.765
    52  = 16.826
    57  = 16.587
    59  = 16.482
    62  = 16.588
    66  = 17.482
    73  = 16.268
    74  = 16.484
    77  = 16.683
    92  = 17.867
    100 = 16.768
}

foreach ($row in $updates.Keys) {
    $ws.Range("E$row").Value = $updates[$row]
}
